# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment Schedule" sheet, shifting "Late" -> O and "Outstanding" -> Q,
# then make "Repayment Schedule" the active sheet with S9 selected
# (previously "Transactions" was the active sheet with B1 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()
$ws.Columns("N:N").Insert()
$ws.Range("S9").Select() | Out-Null
